# "got some simulation results" - fill in the 2-level/3-level SPWM and SVPWM
# capacitor-value comparison tables (without vs with interleaving).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a value, then italicise+shrink a tail run inside it (used for
# the "... with same total load" footnotes that share a cell with the main
# reading).
# ---------------------------------------------------------------------------
function Set-TailNote($addr, $text, $tailStart, $tailLen) {
    $ws.Range($addr).Value = $text
    if ($tailLen -gt 0) {
        $r = $ws.Range($addr).Characters($tailStart, $tailLen)
        $r.Font.Italic = $true
        $r.Font.Size = 8
        $r.Font.Name = "Calibri"
    }
}

# ---------------------------------------------------------------------------
# Helper: build one "section" of the report — a merged bold title row
# followed by a 2-row header/sub-header pair and two data rows (Voltage
# swing / Capacitor value for the same swing), matching the layout already
# used by the first block in the sheet.
# ---------------------------------------------------------------------------
function Build-Section($titleRow, $headerRow, $voltRow, $capRow, $titleText,
                        $withoutVolt, $withVolt, $withVoltTailStart, $withVoltTailLen,
                        $withoutCap, $withCap, $withCapTailStart, $withCapTailLen) {

    $ws.Range("A$titleRow" + ":D$titleRow").Merge()
    $ws.Range("A$titleRow").Value = $titleText
    $tr = $ws.Range("A$titleRow" + ":D$titleRow")
    $tr.Font.Bold = $true
    $tr.Font.Name = "Arial"
    $tr.Font.Size = 12
    $tr.HorizontalAlignment = -4108
    $tr.VerticalAlignment = -4108

    $ws.Range("B$headerRow").Value = "Without-interleaving"
    $ws.Range("C$headerRow").Value = "With-Interleaving (2 inverters using the same dc link)"
    $ws.Range("D$headerRow").Value = "Load Value "
    $hdr = $ws.Range("B$headerRow" + ":C$headerRow")
    $hdr.Font.Bold = $true
    $hdr.Font.Name = "Arial"
    $hdr.Font.Size = 12
    $hdr.HorizontalAlignment = -4108
    $hdr.VerticalAlignment = -4108
    $hdr.WrapText = $true
    $dh = $ws.Range("D$headerRow")
    $dh.Font.Bold = $false
    $dh.Font.Name = "Calibri"
    $dh.Font.Size = 11
    $dh.HorizontalAlignment = -4108
    $dh.VerticalAlignment = -4108

    $ws.Range("A$voltRow").Value = "Voltage swing "
    $av = $ws.Range("A$voltRow")
    $av.Font.Bold = $true
    $av.Font.Name = "Arial"
    $av.Font.Size = 12
    $av.HorizontalAlignment = -4108
    $av.VerticalAlignment = -4108

    Set-TailNote "B$voltRow" $withoutVolt 0 0
    $ws.Range("B$voltRow").HorizontalAlignment = -4108
    $ws.Range("B$voltRow").VerticalAlignment = -4108

    Set-TailNote "C$voltRow" $withVolt $withVoltTailStart $withVoltTailLen
    $cv = $ws.Range("C$voltRow")
    $cv.HorizontalAlignment = -4108
    $cv.VerticalAlignment = -4108
    $cv.WrapText = $true

    $ws.Range("D$voltRow").Value = "8888VA  0.9pf lagging per inverter"
    $dv = $ws.Range("D$voltRow")
    $dv.HorizontalAlignment = -4108
    $dv.VerticalAlignment = -4108

    $ws.Range("A$capRow").Value = "Capacitor value for the same swing (For 1Vpp)"
    $ac = $ws.Range("A$capRow")
    $ac.Font.Bold = $false
    $ac.Font.Name = "Arial"
    $ac.Font.Size = 12
    $ac.HorizontalAlignment = -4108
    $ac.WrapText = $true

    Set-TailNote "B$capRow" $withoutCap 0 0
    $ws.Range("B$capRow").HorizontalAlignment = -4108
    $ws.Range("B$capRow").VerticalAlignment = -4108

    Set-TailNote "C$capRow" $withCap $withCapTailStart $withCapTailLen
    $cc = $ws.Range("C$capRow")
    $cc.HorizontalAlignment = -4108
    $cc.VerticalAlignment = -4108
    $cc.WrapText = $true

    $ws.Range("D$capRow").Value = "8888VA  0.9pf lagging per inverter"
    $dc = $ws.Range("D$capRow")
    $dc.HorizontalAlignment = -4108
    $dc.VerticalAlignment = -4108
}

# ---------------------------------------------------------------------------
# Section 1 (rows 1-4, existing block): re-purposed for "2 level SPWM"
# ---------------------------------------------------------------------------
Build-Section 1 2 3 4 "2 level SPWM (ma = 1, mf = 41,fundf = 50hz)" `
    "269.6 - 269.4 (Vpp = 0.22 V)" `
    "268.8 - 268.5 (Vpp = 0.3 V)                        269.4 - 269.2 (Vpp = 0.15 V) with same total load" 46 55 `
    "100µF  (4 pieces for 2 inverters)" `
    "280µF (2 pieces for 2 inverters)                       80µF (2 pieces for 2 inverters) with same total load" 50 58

# ---------------------------------------------------------------------------
# Section 2 (rows 6-9): "3 level SPWM" — this reuses the figures that used
# to live in the old rows 1/3/4 before the sheet grew.
# ---------------------------------------------------------------------------
Build-Section 6 7 8 9 "3 level SPWM (ma = 1, mf = 41,fundf = 50hz)" `
    "271.1 - 269.7 (Vpp = 1.387V)" `
    "269.6 - 269.3 (Vpp = 0.3V)" 0 0 `
    "1350µF (4 pieces for 2 inverters)" `
    "300µF (2 pieces for 2 inverters)" 0 0

# ---------------------------------------------------------------------------
# Section 3 (rows 11-14): "2 level SVPWM"
# ---------------------------------------------------------------------------
Build-Section 11 12 13 14 "2 level SVPWM (ma = 1, swf = 2050hz,fundf = 50hz)" `
    "269.4 - 268.8 (Vpp = 0.63V)" `
    "268.5 - 267.7 (Vpp = 0.74V)                        269.2 - 268.9 (Vpp = 0.37 V) with same total load" 48 53 `
    "270µF (4 pieces for 2 inverters)" `
    "430µF (2 pieces for 2 inverters)                       60µF (2 pieces for 2 inverters) with same total load" 54 54

# ---------------------------------------------------------------------------
# Section 4 (rows 16-19): "3 level SVPWM"
# ---------------------------------------------------------------------------
Build-Section 16 17 18 19 "3 level SVPWM (ma = 1, swf = 2050hz,fundf = 50hz)" `
    "269.5 - 268.1 (Vpp = 1.47V)" `
    "268.4 - 267.8 (Vpp = 0.6V)" 0 0 `
    "1400µF (4 pieces for 2 inverters)" `
    "250µF (2 pieces for 2 inverters)" 0 0

# ---------------------------------------------------------------------------
# Leftover spacer rows between sections / trailing blank row, same styling
# as the spacer row that always separated the original table from the rest
# of the (previously empty) sheet.
# ---------------------------------------------------------------------------
$ws.Range("E6").VerticalAlignment = -4108
$ws.Range("E7:E9").VerticalAlignment = -4108
$ws.Range("E10").VerticalAlignment = -4108
$ws.Range("E11:E14").VerticalAlignment = -4108
$ws.Range("E15").VerticalAlignment = -4108

$ws.Range("B21").HorizontalAlignment = -4108
$ws.Range("B21").VerticalAlignment = -4108
$ws.Range("C21:D21").VerticalAlignment = -4108

# Drop the leftover columns B:D that used to carry the blank placeholder
# formatting on rows 10 and 15 now that those rows only hold the spacer
# cell in column E.
$ws.Range("B10:D10").Clear()
$ws.Range("B15:D15").Clear()

# ---------------------------------------------------------------------------
# Sheet-level cosmetics: wider columns now that the table has 4 content
# columns instead of 3, new selection/zoom, row heights.
# ---------------------------------------------------------------------------
$ws.Columns("A:B").ColumnWidth = 41.6
$ws.Columns("C:D").ColumnWidth = 41.6

$ws.Rows(1).RowHeight = 29.25
$ws.Rows(6).RowHeight = 34.5
$ws.Rows(7).RowHeight = 30.75
$ws.Rows(8).RowHeight = 30
$ws.Rows(9).RowHeight = 30
$ws.Rows(10).RowHeight = 25.5
$ws.Rows(11).RowHeight = 36.75
$ws.Rows(12).RowHeight = 33.75
$ws.Rows(13).RowHeight = 26.25
$ws.Rows(14).RowHeight = 40.5
$ws.Rows(15).RowHeight = 25.5
$ws.Rows(16).RowHeight = 36.75
$ws.Rows(17).RowHeight = 36.75
$ws.Rows(18).RowHeight = 36.75
$ws.Rows(19).RowHeight = 36.75

$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$ws.Range("B2").Select()
